$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 and 17 coin entries swap positions (WrappedBTC <-> ShibaInu) with updated data

$ws.Range("D2").Value = "29.249.14"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.871.82"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "0.7112"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").Value = "241.65"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "0.3106"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").Value = "0.08380"
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").Value = "1.892.53"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "5.208"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "0.7095"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "91.24"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.000008315"
$ws.Range("E16").Value = "  +6.57%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "29.252.11"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "5.951"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").Value = "242.29"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "2.129.00"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "13.18"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "7.814"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "0.1627"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").Value = "163.35"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").Value = "9.013"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").Value = "18.49"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").Value = "1.502"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").Value = "4.409"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "4.312"
$ws.Range("E31").Value = "  +5.45%  "
$ws.Range("E32").Value = "  -4.46%  "
$ws.Range("D33").Value = "0.05244"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "1.920"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").Value = "0.7475"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("D36").Value = "1.171"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").Value = "2.681"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "0.01856"
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("D39").Value = "2.713"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").Value = "1.151.73"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").Value = "6.363"
$ws.Range("E41").Value = "  +4.08%  "
$ws.Range("D42").Value = "73.09"
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("D43").Value = "0.8847"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("D44").Value = "104.90"
$ws.Range("E44").Value = "  +3.04%  "
$ws.Range("D45").Value = "0.9996"
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("D46").Value = "2.025.06"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").Value = "1.798"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").Value = "0.5189"
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "9.364"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").Value = "0.4294"
$ws.Range("E51").Value = "  +0.97%  "
